$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Who is doing" column (B) to "Tomek" and "Status" column (C) to "Done"
# for rows 2 through 6 (matching existing data validation list options).
$ws.Range("B2:B6").Value = "Tomek"
$ws.Range("C2:C6").Value = "Done"

# Update the active selection to match the edited region.
$ws.Range("C3:C6").Select()
